$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "time_taken", matching the style of the existing header row (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate time_taken column for each data row (F2:F22) with the recorded timestamps
$ws.Range("F2").Value = "2021-10-05 13:39:10.806378"
$ws.Range("F3").Value = "2021-10-05 13:39:10.806390"
$ws.Range("F4").Value = "2021-10-05 13:39:10.806393"
$ws.Range("F5").Value = "2021-10-05 13:39:10.806396"
$ws.Range("F6").Value = "2021-10-05 13:39:10.806399"
$ws.Range("F7").Value = "2021-10-05 13:39:10.806402"
$ws.Range("F8").Value = "2021-10-05 13:39:10.806404"
$ws.Range("F9").Value = "2021-10-05 13:39:10.806407"
$ws.Range("F10").Value = "2021-10-05 13:39:10.806409"
$ws.Range("F11").Value = "2021-10-05 13:39:10.806412"
$ws.Range("F12").Value = "2021-10-05 13:39:10.806414"
$ws.Range("F13").Value = "2021-10-05 13:39:10.806417"
$ws.Range("F14").Value = "2021-10-05 13:39:10.806419"
$ws.Range("F15").Value = "2021-10-05 13:39:10.806421"
$ws.Range("F16").Value = "2021-10-05 13:39:10.806424"
$ws.Range("F17").Value = "2021-10-05 13:39:10.806427"
$ws.Range("F18").Value = "2021-10-05 13:39:10.806429"
$ws.Range("F19").Value = "2021-10-05 13:39:10.806432"
$ws.Range("F20").Value = "2021-10-05 13:39:10.806434"
$ws.Range("F21").Value = "2021-10-05 13:39:10.806437"
$ws.Range("F22").Value = "2021-10-05 13:39:10.806439"
